$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 230-445 get a new "OK" value in column D and a tightened row height,
# except the rows whose C cell already carries the "duplicate/equal pair"
# styling (s=2 / s=3) which are left untouched by this pass.
$skip = @(241, 274, 288, 310, 326, 335, 360, 381, 412, 414, 433)

for ($r = 230; $r -le 445; $r++) {
    if ($skip -contains $r) { continue }
    $ws.Range("D$r").Value = "OK"
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Row 305's word pair gets promoted to bold (matches the existing "bold"
# cell style already used elsewhere in the sheet).
$ws.Range("C305").Font.Bold = $true

# Move the view/selection forward to where work resumed.
$excel.ActiveWindow.ScrollRow = 430
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D440:D445").Select()

Write-Host "done"
